$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.014.54"
$ws.Range("E2").Value = "  -6.81%  "
$ws.Range("D3").Value = "1.413.47"
$ws.Range("E3").Value = "  -7.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.05"
$ws.Range("E6").Value = "  -4.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3657"
$ws.Range("E7").Value = "  -5.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3121"
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.80"
$ws.Range("E9").Value = "  -6.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.039"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06516"
$ws.Range("E11").Value = "  -8.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9996"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.526"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.81"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.203"
$ws.Range("E15").Value = "  -4.83%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001024"
$ws.Range("E16").Value = "  -5.46%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.412.38"
$ws.Range("E17").Value = "  -7.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05691"
$ws.Range("E18").Value = "  -13.81%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.04"
$ws.Range("E20").Value = "  -14.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.635"
$ws.Range("E21").Value = "  -7.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.78"
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.262"
$ws.Range("D25").Value = "20.001.10"
$ws.Range("E25").Value = "  -6.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.282"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "133.44"
$ws.Range("E27").Value = "  -10.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.22"
$ws.Range("E28").Value = "  -5.91%  "
$ws.Range("D29").Value = "1.570.72"
$ws.Range("E29").Value = "  -7.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.77"
$ws.Range("E30").Value = "  -5.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.957"
$ws.Range("E31").Value = "  -17.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.308"
$ws.Range("E32").Value = "  -11.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8251"
$ws.Range("E33").Value = "  -13.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07699"
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.418"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.483"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05920"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.923"
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.0000"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02086"
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.54"
$ws.Range("E41").Value = "  -5.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1910"
$ws.Range("E42").Value = "  -5.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.097"
$ws.Range("E43").Value = "  -6.55%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5333"
$ws.Range("E44").Value = "  -6.88%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.30"
$ws.Range("E45").Value = "  -6.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.537"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5222"
$ws.Range("E47").Value = "  -5.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "115.01"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.779"
$ws.Range("E49").Value = "  -5.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.039"
$ws.Range("E50").Value = "  -10.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  +0.02%  "
